$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

# New "Date" timestamps for column B (rows 2-16)
$dates = @(
    "Mon Jan 09 17:50:24 EST 2023",
    "Mon Jan 09 17:50:33 EST 2023",
    "Mon Jan 09 17:50:42 EST 2023",
    "Mon Jan 09 17:50:52 EST 2023",
    "Mon Jan 09 17:51:02 EST 2023",
    "Mon Jan 09 17:51:11 EST 2023",
    "Mon Jan 09 17:51:21 EST 2023",
    "Mon Jan 09 17:51:31 EST 2023",
    "Mon Jan 09 17:51:41 EST 2023",
    "Mon Jan 09 17:51:51 EST 2023",
    "Mon Jan 09 17:52:01 EST 2023",
    "Mon Jan 09 17:52:11 EST 2023",
    "Mon Jan 09 17:52:21 EST 2023",
    "Mon Jan 09 17:52:30 EST 2023",
    "Mon Jan 09 17:52:40 EST 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $dates[$i]
    $ws.Range("C$row").Value = "Y"
}

# Touching H1 and clearing it stretches the worksheet's used-range/dimension
# to column H (matching the saved workbook's recorded extent) without
# leaving any visible content behind.
$ws.Range("H1").Value = "temp"
$ws.Range("H1").ClearContents()

$ws.Range("C2:C16").Select() | Out-Null
